$d = $word.ActiveDocument

# 1. Title text (appears twice: Heading1 and bold run near the bottom)
$d.Content.Find.Execute(
    "Play Mercy of the Gods for Free - Review and Gameplay Mechanics", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Mercy of the Gods Free - Exciting Ancient Egyptian Slot", 2)

# 2. "What we like" bullet list items
$d.Content.Find.Execute(
    "Immersive ancient Egypt theme", $true, $false, $false, $false, $false,
    $true, 1, $false, "Immersive Ancient Egypt theme", 2)

$d.Content.Find.Execute(
    "High volatility with significant payouts", $true, $false, $false, $false, $false,
    $true, 1, $false, "High volatility and payouts for big wins", 2)

$d.Content.Find.Execute(
    "Wild symbol increases player's chances of earning large payouts", $true, $false, $false, $false, $false,
    $true, 1, $false, "Exciting gameplay mechanics", 2)

$d.Content.Find.Execute(
    "Relaxing background music", $true, $false, $false, $false, $false,
    $true, 1, $false, "Wild symbol increases chances of big payouts", 2)

# 3. "What we don't like" bullet list items
$d.Content.Find.Execute(
    "Significant plays before a win due to high volatility", $true, $false, $false, $false, $false,
    $true, 1, $false, "Potentially long runs without a win", 2)

$d.Content.Find.Execute(
    "No bonus game or jackpot feature", $true, $false, $false, $false, $false,
    $true, 1, $false, "Limited number of symbols on the game grid", 2)

# 4. Meta description (italic text)
$d.Content.Find.Execute(
    "Explore the immersive ancient Egypt theme and exciting gameplay mechanics for free by playing Mercy of the Gods. Review provides pros and cons.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Read this review of Mercy of the Gods and play it free. Experience thrilling gameplay and big wins in this Ancient Egyptian themed slot.", 2)
